$wb = $excel.ActiveWorkbook

# ALC row 19 (Leve Item ID 7015)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3337.889
$ws.Range("I19").Value = 2398.3333
$ws.Range("J19").Value = 3807.6667
$ws.Range("K19").Value = 2398.3333
$ws.Range("L19").Value = 3807.6667
$ws.Range("M19").Value = -2223.3333
$ws.Range("N19").Value = -4157.6667

# ALC row 34 (Leve Item ID 2160)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3168.4
$ws.Range("I34").Value = 3168.4
$ws.Range("K34").Value = 3168.4
$ws.Range("M34").Value = -2965.4

# ALC row 36 (Leve Item ID 2160)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 3168.4
$ws.Range("I36").Value = 3168.4
$ws.Range("K36").Value = 3168.4
$ws.Range("M36").Value = -2453.4

# ALC row 40 (Leve Item ID 5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3595.2058
$ws.Range("I40").Value = 3247.389
$ws.Range("J40").Value = 3986.5
$ws.Range("K40").Value = 3247.389
$ws.Range("L40").Value = 3986.5
$ws.Range("M40").Value = -3072.389
$ws.Range("N40").Value = -4336.5

# ALC row 80 (Leve Item ID 12605)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1313.7059
$ws.Range("I80").Value = 1336.909
$ws.Range("J80").Value = 1271.1666
$ws.Range("K80").Value = 4010.727
$ws.Range("L80").Value = 3813.4998
$ws.Range("M80").Value = -3012.727
$ws.Range("N80").Value = -5809.4998

# ALC row 83 (Leve Item ID 12605)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1313.7059
$ws.Range("I83").Value = 1336.909
$ws.Range("J83").Value = 1271.1666
$ws.Range("K83").Value = 12032.181
$ws.Range("L83").Value = 11440.4994
$ws.Range("M83").Value = -7040.181
$ws.Range("N83").Value = -21424.4994

# ALC row 133 (Leve Item ID 41856)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 74250
$ws.Range("I133").Value = 55000
$ws.Range("J133").Value = 93500
$ws.Range("K133").Value = 55000
$ws.Range("L133").Value = 93500
$ws.Range("M133").Value = -49940
$ws.Range("N133").Value = -103620

# ALC row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3469.0303
$ws.Range("I137").Value = 2372.95
$ws.Range("K137").Value = 7118.849999999999
$ws.Range("M137").Value = -4568.849999999999

# ALC row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3313.7026
$ws.Range("J138").Value = 3859.1304
$ws.Range("L138").Value = 11577.3912
$ws.Range("N138").Value = -21857.3912

# ARM row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2743.6667
$ws.Range("I32").Value = 2743.6667
$ws.Range("K32").Value = 2743.6667
$ws.Range("M32").Value = -2456.6667

# ARM row 45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12649
$ws.Range("I45").Value = 14170.857
$ws.Range("K45").Value = 14170.857
$ws.Range("M45").Value = -13793.857

# ARM row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 277261.78
$ws.Range("I74").Value = 439751.97
$ws.Range("J74").Value = 10313.643
$ws.Range("K74").Value = 439751.97
$ws.Range("L74").Value = 10313.643
$ws.Range("M74").Value = -438877.97
$ws.Range("N74").Value = -12061.643

# ARM row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 277261.78
$ws.Range("I77").Value = 439751.97
$ws.Range("J77").Value = 10313.643
$ws.Range("K77").Value = 2198759.85
$ws.Range("L77").Value = 51568.215
$ws.Range("M77").Value = -2194391.85
$ws.Range("N77").Value = -60304.215

# ARM row 122 (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3940.5
$ws.Range("I122").Value = 3711.889
$ws.Range("J122").Value = 5998
$ws.Range("K122").Value = 11135.667
$ws.Range("L122").Value = 17994
$ws.Range("M122").Value = -8685.667000000001
$ws.Range("N122").Value = -22894

# ARM row 135 (Leve Item ID 42016)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 116781.25
$ws.Range("J135").Value = 116781.25
$ws.Range("L135").Value = 116781.25
$ws.Range("N135").Value = -126921.25

# BSM row 108 (Leve Item ID 25643)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 87375.8
$ws.Range("J108").Value = 87375.8
$ws.Range("L108").Value = 87375.8
$ws.Range("N108").Value = -95055.8

# CRP row 10 (Leve Item ID 1997)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 309.26666
$ws.Range("I10").Value = 149.25
$ws.Range("J10").Value = 492.14285
$ws.Range("K10").Value = 149.25
$ws.Range("L10").Value = 492.14285
$ws.Range("M10").Value = -10.25
$ws.Range("N10").Value = -770.14285

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3255.3142
$ws.Range("I31").Value = 1292.8572
$ws.Range("J31").Value = 4563.619
$ws.Range("K31").Value = 1292.8572
$ws.Range("L31").Value = 4563.619
$ws.Range("M31").Value = -997.8571999999999
$ws.Range("N31").Value = -5153.619

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3255.3142
$ws.Range("I34").Value = 1292.8572
$ws.Range("J34").Value = 4563.619
$ws.Range("K34").Value = 1292.8572
$ws.Range("L34").Value = 4563.619
$ws.Range("M34").Value = -1090.8572
$ws.Range("N34").Value = -4967.619

# CRP row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3506.5151
$ws.Range("I58").Value = 1371.909
$ws.Range("K58").Value = 1371.909
$ws.Range("M58").Value = -1168.909

# CRP row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 337376.44
$ws.Range("J132").Value = 308847
$ws.Range("L132").Value = 926541
$ws.Range("N132").Value = -931601

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4421.5864
$ws.Range("I134").Value = 3505.25
$ws.Range("K134").Value = 10515.75
$ws.Range("M134").Value = -7980.75

# CRP row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3506.5151
$ws.Range("I136").Value = 1371.909
$ws.Range("K136").Value = 4115.727000000001
$ws.Range("M136").Value = -1565.727000000001

# CUL row 2 (Leve Item ID 4847)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3232
$ws.Range("J2").Value = 3926.111
$ws.Range("L2").Value = 23556.666
$ws.Range("N2").Value = -23782.666

# CUL row 4 (Leve Item ID 4650)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4207308.5
$ws.Range("I4").Value = 3219830
$ws.Range("K4").Value = 9659490
$ws.Range("M4").Value = -9659378

# CUL row 23 (Leve Item ID 4858)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 166917.17
$ws.Range("I23").Value = 333468.66
$ws.Range("J23").Value = 365.66666
$ws.Range("K23").Value = 1000405.98
$ws.Range("L23").Value = 1096.99998
$ws.Range("M23").Value = -1000170.98
$ws.Range("N23").Value = -1566.99998

# CUL row 31 (Leve Item ID 4710)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 548.3333
$ws.Range("I31").Value = 548.3333
$ws.Range("K31").Value = 1644.9999
$ws.Range("M31").Value = -1356.9999

# CUL row 34 (Leve Item ID 4749)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1823.125
$ws.Range("I34").Value = 114.666664
$ws.Range("J34").Value = 2067.1904
$ws.Range("K34").Value = 343.999992
$ws.Range("L34").Value = 6201.5712
$ws.Range("M34").Value = -259.999992
$ws.Range("N34").Value = -6369.5712

# CUL row 38 (Leve Item ID 4860)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 42.272728
$ws.Range("I38").Value = 30.307692
$ws.Range("J38").Value = 59.555557
$ws.Range("K38").Value = 90.92307599999999
$ws.Range("L38").Value = 178.666671
$ws.Range("M38").Value = 256.076924
$ws.Range("N38").Value = -872.666671

# CUL row 39 (Leve Item ID 4712)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6299.6
$ws.Range("J39").Value = 8832.666999999999
$ws.Range("L39").Value = 26498.001
$ws.Range("N39").Value = -27086.001

# CUL row 55 (Leve Item ID 4733)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1393.3334
$ws.Range("J55").Value = 1512.1
$ws.Range("L55").Value = 4536.299999999999
$ws.Range("N55").Value = -4890.299999999999

# CUL row 63 (Leve Item ID 12866)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 499.5
$ws.Range("I63").Value = 499.5
$ws.Range("K63").Value = 1498.5
$ws.Range("M63").Value = -749.5

# CUL row 66 (Leve Item ID 12866)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 499.5
$ws.Range("I66").Value = 499.5
$ws.Range("K66").Value = 4495.5
$ws.Range("M66").Value = -751.5

# CUL row 114 (Leve Item ID 27865)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2227.6667
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# CUL row 116 (Leve Item ID 27866)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1763.25
$ws.Range("I116").Value = 1763.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5289.75
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1847.75
$ws.Range("N116").ClearContents()

# CUL row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11909016
$ws.Range("I131").Value = 41667544
$ws.Range("J131").Value = 5604.25
$ws.Range("K131").Value = 125002632
$ws.Range("L131").Value = 16812.75
$ws.Range("M131").Value = -124997592
$ws.Range("N131").Value = -26892.75

# GSM row 48 (Leve Item ID 4337)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 60000
$ws.Range("J48").Value = 60000
$ws.Range("L48").Value = 60000
$ws.Range("N48").Value = -60970

# GSM row 130 (Leve Item ID 34692)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 82000
$ws.Range("J130").Value = 82000
$ws.Range("L130").Value = 82000
$ws.Range("N130").Value = -92040

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9039.111000000001
$ws.Range("J132").Value = 13166.333
$ws.Range("L132").Value = 39498.999
$ws.Range("N132").Value = -44558.999

# LTW row 20 (Leve Item ID 4308)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15000000
$ws.Range("I20").Value = 15000000
$ws.Range("K20").Value = 15000000
$ws.Range("M20").Value = -14999774

# LTW row 42 (Leve Item ID 4333)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 50000000
$ws.Range("I42").Value = 50000000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 50000000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -49999437
$ws.Range("N42").ClearContents()

# LTW row 49 (Leve Item ID 4333)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 50000000
$ws.Range("I49").Value = 50000000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 50000000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -49999853
$ws.Range("N49").ClearContents()

# LTW row 87 (Leve Item ID 10926)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 65000
$ws.Range("J87").Value = 65000
$ws.Range("L87").Value = 65000
$ws.Range("N87").Value = -67246

# LTW row 90 (Leve Item ID 10926)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 65000
$ws.Range("J90").Value = 65000
$ws.Range("L90").Value = 195000
$ws.Range("N90").Value = -206232

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4942.9644
$ws.Range("I136").Value = 3137.65
$ws.Range("K136").Value = 9412.950000000001
$ws.Range("M136").Value = -6862.950000000001

# WVR row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3108.9788
$ws.Range("I132").Value = 2138.2974
$ws.Range("K132").Value = 6414.8922
$ws.Range("M132").Value = -3884.8922

# WVR row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2536.3333
$ws.Range("I136").Value = 1064.3914
$ws.Range("K136").Value = 3193.1742
$ws.Range("M136").Value = -643.1741999999999
